$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new requirement R4 (row 6): "Destruir edificio"
$ws.Range("E6").Value = "Destruir edificio"
$ws.Range("F6").Value = "Cuando el personaje golpea al edificio, este empieza a destruirse hasta derrumbarse."

# Update the "Entradas" cell for R1: "Evento de elección" -> "Objeto de tipo jugador"
$ws.Range("G3").Value = "Objeto de tipo jugador"

$ws.Range("G6").Value = "Coordenadas de ambos objetos que entran en colisión"
$ws.Range("H6").Value = "El edificio va desaparece en las coordenadas de colisión."

# Fill in the new requirement R5 (row 7): "Entrar en bonus"
$ws.Range("E7").Value = "Entrar en bonus"
$ws.Range("F7").Value = "Cuando el personaje entra a cierta parte de la escena, cae en un nivel bonus, entrando en otra escena."
$ws.Range("G7").Value = "Coordenadas del objeto y coordenadas a comparar"
$ws.Range("H7").Value = "Se cambia la escena, cargando la escena del bonus respectivo."

# Rows 6 and 7 now wrap like the other filled-in requirement rows
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2

# Column G needs to widen to fit the new text
$ws.Columns.Item(7).ColumnWidth = 27

# Move the active selection to G7, matching the end-state cursor position
$ws.Range("G7").Select()
